$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 6 (data trimmed to 4 data rows instead of 5)
$ws.Rows.Item(6).Delete()

# Update data rows 2-5 with new values (A:AH)
$ws.Cells.Item(2,1).Value = 45070.50694444445
$ws.Cells.Item(2,2).Value = 15.374
$ws.Cells.Item(2,3).Value = 10.112
$ws.Cells.Item(2,4).Value = 3.717
$ws.Cells.Item(2,5).Value = 32.732
$ws.Cells.Item(2,6).Value = 25.126
$ws.Cells.Item(2,7).Value = 11.943
$ws.Cells.Item(2,8).Value = 36.363
$ws.Cells.Item(2,9).Value = 18.615
$ws.Cells.Item(2,10).Value = 7.558
$ws.Cells.Item(2,11).Value = 11.183
$ws.Cells.Item(2,12).Value = 12.932
$ws.Cells.Item(2,13).Value = 13.607
$ws.Cells.Item(2,14).Value = 3.86
$ws.Cells.Item(2,15).Value = 12.031
$ws.Cells.Item(2,16).Value = 16.608
$ws.Cells.Item(2,17).Value = 10.57
$ws.Cells.Item(2,18).Value = 3.118
$ws.Cells.Item(2,19).Value = 1.764
$ws.Cells.Item(2,20).Value = 175.75
$ws.Cells.Item(2,21).Value = 33.355
$ws.Cells.Item(2,22).Value = 11.105
$ws.Cells.Item(2,23).Value = 21.604
$ws.Cells.Item(2,24).Value = 11.537
$ws.Cells.Item(2,25).Value = 2.926
$ws.Cells.Item(2,26).Value = 18.421
$ws.Cells.Item(2,27).Value = 9.809
$ws.Cells.Item(2,28).Value = 8.874
$ws.Cells.Item(2,29).Value = 10.611
$ws.Cells.Item(2,30).Value = 13.63
$ws.Cells.Item(2,31).Value = 3.312
$ws.Cells.Item(2,32).Value = 32.661
$ws.Cells.Item(2,33).Value = 5.894
$ws.Cells.Item(2,34).Value = 13.883
$ws.Cells.Item(3,1).Value = 45070.51388888889
$ws.Cells.Item(3,2).Value = 12.011
$ws.Cells.Item(3,3).Value = 8.333
$ws.Cells.Item(3,4).Value = 1.629
$ws.Cells.Item(3,5).Value = 25.992
$ws.Cells.Item(3,6).Value = 20.404
$ws.Cells.Item(3,7).Value = 9.326
$ws.Cells.Item(3,8).Value = 36.495
$ws.Cells.Item(3,9).Value = 14.543
$ws.Cells.Item(3,10).Value = 6.187
$ws.Cells.Item(3,11).Value = 8.978
$ws.Cells.Item(3,12).Value = 10.371
$ws.Cells.Item(3,13).Value = 10.96
$ws.Cells.Item(3,14).Value = 3.02
$ws.Cells.Item(3,15).Value = 9.399
$ws.Cells.Item(3,16).Value = 13.16
$ws.Cells.Item(3,17).Value = 8.267
$ws.Cells.Item(3,18).Value = 1.435
$ws.Cells.Item(3,19).Value = 0.892
$ws.Cells.Item(3,20).Value = 135.743
$ws.Cells.Item(3,21).Value = 26.329
$ws.Cells.Item(3,22).Value = 8.676
$ws.Cells.Item(3,23).Value = 17.272
$ws.Cells.Item(3,24).Value = 9.252
$ws.Cells.Item(3,25).Value = 1.875
$ws.Cells.Item(3,26).Value = 17.584
$ws.Cells.Item(3,27).Value = 7.663
$ws.Cells.Item(3,28).Value = 6.954
$ws.Cells.Item(3,29).Value = 8.219
$ws.Cells.Item(3,30).Value = 10.906
$ws.Cells.Item(3,31).Value = 1.246
$ws.Cells.Item(3,32).Value = 33.585
$ws.Cells.Item(3,33).Value = 4.683
$ws.Cells.Item(3,34).Value = 10.847
$ws.Cells.Item(4,1).Value = 45070.52083333334
$ws.Cells.Item(4,2).Value = 0.382
$ws.Cells.Item(4,3).Value = 0.17
$ws.Cells.Item(4,4).Value = 0.773
$ws.Cells.Item(4,5).Value = 0.847
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = 0
$ws.Cells.Item(4,8).Value = 7.849
$ws.Cells.Item(4,9).Value = 0.582
$ws.Cells.Item(4,10).Value = 0.126
$ws.Cells.Item(4,11).Value = 0.264
$ws.Cells.Item(4,12).Value = 0.161
$ws.Cells.Item(4,13).Value = 0.008
$ws.Cells.Item(4,14).Value = 0
$ws.Cells.Item(4,15).Value = 0.376
$ws.Cells.Item(4,16).Value = 0.595
$ws.Cells.Item(4,17).Value = 0.591
$ws.Cells.Item(4,18).Value = 0.873
$ws.Cells.Item(4,19).Value = 0.338
$ws.Cells.Item(4,20).Value = 0
$ws.Cells.Item(4,21).Value = 1.63
$ws.Cells.Item(4,22).Value = 0.347
$ws.Cells.Item(4,23).Value = 0.985
$ws.Cells.Item(4,24).Value = 0.546
$ws.Cells.Item(4,25).Value = 0.477
$ws.Cells.Item(4,26).Value = 3.235
$ws.Cells.Item(4,27).Value = 0.293
$ws.Cells.Item(4,28).Value = 0.387
$ws.Cells.Item(4,29).Value = 0.446
$ws.Cells.Item(4,30).Value = 0.382
$ws.Cells.Item(4,31).Value = 0.715
$ws.Cells.Item(4,32).Value = 7.953
$ws.Cells.Item(4,33).Value = 0.044
$ws.Cells.Item(4,34).Value = 0.446
$ws.Cells.Item(5,1).Value = 45070.52777777778
$ws.Cells.Item(5,2).Value = 0.88
$ws.Cells.Item(5,3).Value = 0.57
$ws.Cells.Item(5,4).Value = 0.58
$ws.Cells.Item(5,5).Value = 1.96
$ws.Cells.Item(5,6).Value = 1.07
$ws.Cells.Item(5,7).Value = 0.92
$ws.Cells.Item(5,8).Value = 4.56
$ws.Cells.Item(5,9).Value = 1.16
$ws.Cells.Item(5,10).Value = 0.25
$ws.Cells.Item(5,11).Value = 0.67
$ws.Cells.Item(5,12).Value = 0.64
$ws.Cells.Item(5,13).Value = 0.6
$ws.Cells.Item(5,14).Value = 0.14
$ws.Cells.Item(5,15).Value = 0.75
$ws.Cells.Item(5,16).Value = 0.92
$ws.Cells.Item(5,17).Value = 0.86
$ws.Cells.Item(5,18).Value = 0.66
$ws.Cells.Item(5,19).Value = 0.25
$ws.Cells.Item(5,20).Value = 4.54
$ws.Cells.Item(5,21).Value = 2.15
$ws.Cells.Item(5,22).Value = 0.69
$ws.Cells.Item(5,23).Value = 1.13
$ws.Cells.Item(5,24).Value = 0.71
$ws.Cells.Item(5,25).Value = 0.42
$ws.Cells.Item(5,26).Value = 1.79
$ws.Cells.Item(5,27).Value = 0.6
$ws.Cells.Item(5,28).Value = 0.63
$ws.Cells.Item(5,29).Value = 0.74
$ws.Cells.Item(5,30).Value = 0.85
$ws.Cells.Item(5,31).Value = 0.52
$ws.Cells.Item(5,32).Value = 4.45
$ws.Cells.Item(5,33).Value = 0.28
$ws.Cells.Item(5,34).Value = 0.85

# Adjust column widths to match new auto-fit values
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667
$ws.Columns.Item(6).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667
$ws.Columns.Item(9).ColumnWidth = 7.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667
$ws.Columns.Item(14).ColumnWidth = 5.166666666666667
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666
$ws.Columns.Item(21).ColumnWidth = 7.166666666666667
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667
$ws.Columns.Item(23).ColumnWidth = 7.166666666666667
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667
$ws.Columns.Item(26).ColumnWidth = 7.166666666666667
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667
